# Add a "canonical SMILES" column (D) to the microstate list sheet.
# For every microstate row this is identical to the existing "canonical
# isomeric SMILES" (column C) EXCEPT for SM16_micro005, whose canonical
# (non-isomeric) SMILES drops the cis/trans bond-stereo markers ("/","\").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -----------------------------------------------------
$ws.Range("D2").Value = "canonical SMILES"

# --- data rows --------------------------------------------------------
# row -> canonical (non-isomeric) SMILES value for column D
$canonical = @{
    3  = "c1cc(c(c(c1)Cl)C(=[OH+])N=c2cc[nH]cc2)Cl"
    4  = "c1cc(c(c(c1)Cl)C(=O)Nc2ccncc2)Cl"
    5  = "c1cc(c(c(c1)Cl)C(=[OH+])[N-]c2ccncc2)Cl"
    6  = "c1cc(c(c(c1)Cl)C(=[NH+]c2cc[nH+]cc2)[O-])Cl"
    7  = "c1cc(c(c(c1)Cl)C(=Nc2ccncc2)[O-])Cl"
    8  = "c1cc(c(c(c1)Cl)C(=[OH+])Nc2ccncc2)Cl"
    9  = "c1cc(c(c(c1)Cl)C(=O)N=c2cc[nH]cc2)Cl"
    10 = "c1cc(c(c(c1)Cl)C(=[OH+])Nc2cc[nH+]cc2)Cl"
}

foreach ($r in $canonical.Keys) {
    $ws.Range("D$r").Value = $canonical[$r]
}

# --- column width for the new column D -------------------------------
$ws.Range("D1").ColumnWidth = 36
